$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------------------
$ws.Range("F5").Value = "Document review"
$ws.Range("G5").Value = "Design"
$ws.Range("H5").Value = "Normal"
$ws.Range("I5").Value = "Detailed design"
$ws.Range("J5").Value = "Medium"
$ws.Range("K5").Value = "SDD_MediaManager.xls, revesion 128"
$ws.Range("L5").Value = "Design issue"
$ws.Range("M5").Value = "TuTM"
$ws.Range("N5").Value = "BaoNDD"
$ws.Range("O5").Value = 40838
$ws.Range("O5").NumberFormat = "mm-dd-yy"
$ws.Range("P5").Value = "Phần này không cần làm"
$ws.Range("Q5").Value = "Không cần quá chi tiết như vậy"
$ws.Range("S5").Value = 40838
$ws.Range("T5").Value = 40838
$ws.Rows("5").RowHeight = 45

# ---------------------------------------------------------------------------
# Row 6
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = "Document review"
$ws.Range("G6").Value = "Design"
$ws.Range("H6").Value = "Normal"
$ws.Range("I6").Value = "Detailed design"
$ws.Range("J6").Value = "Medium"
$ws.Range("K6").Value = "SDD_MediaManager.xls, revesion 128"
$ws.Range("L6").Value = "Design issue"
$ws.Range("M6").Value = "TuTM"
$ws.Range("N6").Value = "BaoNDD"
$ws.Range("O6").Value = 40838
$ws.Range("O6").NumberFormat = "mm-dd-yy"
$ws.Range("P6").Value = "Phần này không cần làm"
$ws.Range("Q6").Value = "Không cần quá chi tiết như vậy"
$ws.Range("S6").Value = 40838
$ws.Range("T6").Value = 40838
$ws.Rows("6").RowHeight = 45

# ---------------------------------------------------------------------------
# Row 7
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = "Document review"
$ws.Range("G7").Value = "Design"
$ws.Range("H7").Value = "Normal"
$ws.Range("I7").Value = "Detailed design"
$ws.Range("J7").Value = "Medium"
$ws.Range("K7").Value = "SDD_ContentManager_New&Edit.xls, revesion 128"
$ws.Range("L7").Value = "Design issue"
$ws.Range("M7").Value = "BaoNDD"
$ws.Range("N7").Value = "TuTM"
$ws.Range("O7").Value = 40838
$ws.Range("O7").NumberFormat = "mm-dd-yy"
$ws.Range("P7").Value = "Sửa lại giá trị tương ứng cho đồng nhất"
$ws.Range("Q7").Value = "Sai sót của developer"
$ws.Range("S7").Value = 40838
$ws.Range("T7").Value = 40838
$ws.Rows("7").RowHeight = 30

# ---------------------------------------------------------------------------
# Row 8
# ---------------------------------------------------------------------------
$ws.Range("F8").Value = "Document review"
$ws.Range("G8").Value = "Design"
$ws.Range("H8").Value = "Normal"
$ws.Range("I8").Value = "Detailed design"
$ws.Range("J8").Value = "Medium"
$ws.Range("K8").Value = "SDD_ContentManager_New&Edit.xls, revesion 128"
$ws.Range("L8").Value = "Design issue"
$ws.Range("M8").Value = "BaoNDD"
$ws.Range("N8").Value = "TuTM"
$ws.Range("O8").Value = 40838
$ws.Range("O8").NumberFormat = "mm-dd-yy"
$ws.Range("P8").Value = "Thêm catid"
$ws.Range("Q8").Value = "Sai sót của developer"
$ws.Range("S8").Value = 40838
$ws.Range("T8").Value = 40838
$ws.Rows("8").RowHeight = 30

# ---------------------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------------------
$ws.Range("F9").Value = "Document review"
$ws.Range("G9").Value = "Design"
$ws.Range("H9").Value = "Normal"
$ws.Range("I9").Value = "Detailed design"
$ws.Range("J9").Value = "Medium"
$ws.Range("K9").Value = "SDD_ContentManager_New&Edit.xls, revesion 128"
$ws.Range("L9").Value = "Design issue"
$ws.Range("M9").Value = "TuanVA"
$ws.Range("N9").Value = "TuTM"
$ws.Range("O9").Value = 40838
$ws.Range("O9").NumberFormat = "mm-dd-yy"
$ws.Range("P9").Value = "Thêm event select change của combox section"
$ws.Range("Q9").Value = "Sai sót của developer"
$ws.Range("S9").Value = 40838
$ws.Range("T9").Value = 40838
$ws.Rows("9").RowHeight = 30

# ---------------------------------------------------------------------------
# Row 10
# ---------------------------------------------------------------------------
$ws.Range("F10").Value = "Document review"
$ws.Range("G10").Value = "Design"
$ws.Range("H10").Value = "Normal"
$ws.Range("I10").Value = "Detailed design"
$ws.Range("J10").Value = "Medium"
$ws.Range("K10").Value = "SDD_MediaManager.xls, revesion 128"
$ws.Range("L10").Value = "Design issue"
$ws.Range("M10").Value = "TuTM"
$ws.Range("N10").Value = "TuanVA"
$ws.Range("O10").Value = 40838
$ws.Range("O10").NumberFormat = "mm-dd-yy"
$ws.Range("P10").Value = "Dùng button hoặc dùng link"
$ws.Range("Q10").Value = "Sai sót của developer"
$ws.Range("S10").Value = 40838
$ws.Range("T10").Value = 40838
$ws.Rows("10").RowHeight = 30

# ---------------------------------------------------------------------------
# Row 11
# ---------------------------------------------------------------------------
$ws.Range("F11").Value = "Document review"
$ws.Range("G11").Value = "Design"
$ws.Range("H11").Value = "Normal"
$ws.Range("I11").Value = "Detailed design"
$ws.Range("J11").Value = "Medium"
$ws.Range("K11").Value = "SDD_MediaManager.xls, revesion 128"
$ws.Range("L11").Value = "Design issue"
$ws.Range("M11").Value = "TuanVA"
$ws.Range("N11").Value = "TuanVA"
$ws.Range("O11").Value = 40838
$ws.Range("O11").NumberFormat = "mm-dd-yy"
$ws.Range("P11").Value = "Bỏ giá trị này"
$ws.Range("Q11").Value = "Sai sót của developer"
$ws.Range("S11").Value = 40838
$ws.Range("T11").Value = 40838
$ws.Rows("11").RowHeight = 30

# ---------------------------------------------------------------------------
# Row 12
# ---------------------------------------------------------------------------
$ws.Range("F12").Value = "Document review"
$ws.Range("G12").Value = "Design"
$ws.Range("H12").Value = "Normal"
$ws.Range("I12").Value = "Detailed design"
$ws.Range("J12").Value = "Medium"
$ws.Range("K12").Value = "SDD_MediaManager.xls, revesion 128"
$ws.Range("L12").Value = "Design issue"
$ws.Range("M12").Value = "TuanVA"
$ws.Range("N12").Value = "TuanVA"
$ws.Range("O12").Value = 40838
$ws.Range("O12").NumberFormat = "mm-dd-yy"
$ws.Range("P12").Value = "Bỏ giá trị này"
$ws.Range("Q12").Value = "Sai sót của developer"
$ws.Range("S12").Value = 40838
$ws.Range("T12").Value = 40838
$ws.Rows("12").RowHeight = 30

# ---------------------------------------------------------------------------
# Row 13
# ---------------------------------------------------------------------------
$ws.Range("F13").Value = "Document review"
$ws.Range("G13").Value = "Design"
$ws.Range("H13").Value = "Normal"
$ws.Range("I13").Value = "Detailed design"
$ws.Range("J13").Value = "Medium"
$ws.Range("K13").Value = "SDD_MediaManager.xls, revesion 128"
$ws.Range("L13").Value = "Design issue"
$ws.Range("M13").Value = "BaoNDD"
$ws.Range("N13").Value = "TuanVA"
$ws.Range("O13").Value = 40838
$ws.Range("O13").NumberFormat = "mm-dd-yy"
$ws.Range("P13").Value = "Thêm event select chọn trong folder"
$ws.Range("Q13").Value = "Sai sót của developer"
$ws.Range("S13").Value = 40838
$ws.Range("T13").Value = 40838
$ws.Rows("13").RowHeight = 30

# ---------------------------------------------------------------------------
# Row 14 / 15 height adjustments
# ---------------------------------------------------------------------------
$ws.Rows("14").RowHeight = 36
$ws.Rows("15").RowHeight = 60

# ---------------------------------------------------------------------------
# Column width tweaks (K / O) - closest achievable approximation
# ---------------------------------------------------------------------------
$ws.Columns("K").ColumnWidth = 33.25
$ws.Columns("O").ColumnWidth = 9.75

# ---------------------------------------------------------------------------
# Sheet view: scroll position + active selection
# ---------------------------------------------------------------------------
$ws.Range("R1").Select()
$excel.ActiveWindow.ScrollColumn = 10
